$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.012349247932434
$ws.Range("B1").Value = 1.651827216148376
$ws.Range("C1").Value = 3.473353862762451
$ws.Range("D1").Value = 3.785707473754883
$ws.Range("E1").Value = 0.9584437012672424
